$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1.445647641019636;   C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987; F = 0; G = 6.82939032824165 }
    3 = @{ B = 0.04172184405617529; C = 0.3048912486333797;  D = 189.6080260415259;   E = 0.5333859586016987; F = 1; G = 190.4880250928172 }
    4 = @{ B = 0.04172184405617529; C = 1.626987699542094;  D = 0.1496068669990043;  E = 0.5333859586016987; F = 1; G = 2.351702369198972 }
    5 = @{ B = 3.272327238179451;   C = 1.626987699542094;  D = 3.223369029078222;   E = 0.5333859586016987; F = 1; G = 8.656069925401464 }
    6 = @{ B = 0.01253208636536152; C = 41249014.21622031;  D = 3.223369029078222;   E = 2797.565817734744;  F = 1; G = 41251815.01793917 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
